$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: man_porpuesta -> man_propuesta (PROPUESTA section, field name)
$ws.Range("B54").Value = "man_propuesta"

# Fix field name: id_servicio -> id_servicios (SERVICIOS section, field name)
$ws.Range("B63").Value = "id_servicios"

# tpDoc_mujer size correction: 5 -> 30
$ws.Range("C25").Value = 30

# fkDoc_mujer size/type correction: 11/INT -> 19/BIGINT
$ws.Range("C74").Value = 19
$ws.Range("D74").Value = "BIGINT"
